$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -0.7562163359180574
$ws.Range("D2").Value = 0.4575413224477005

$ws.Range("C3").Value = -0.5735825733982975
$ws.Range("D3").Value = 0.5720678861681465

$ws.Range("C4").Value = 1.665665569537948
$ws.Range("D4").Value = 0.1099573845596935

$ws.Range("C5").Value = 0.07426082270509883
$ws.Range("D5").Value = 0.9414739476839102

$ws.Range("C6").Value = 0.2524721524243098
$ws.Range("D6").Value = 0.8030200330154873

$ws.Range("C7").Value = 2.340329406112874
$ws.Range("D7").Value = 0.02874110344570058

$ws.Range("C8").Value = 0.7047094043381656
$ws.Range("D8").Value = 0.4883850419611599

$ws.Range("C9").Value = 2.611529893878576
$ws.Range("D9").Value = 0.01593304077264412

$ws.Range("C10").Value = 0.867240756290948
$ws.Range("D10").Value = 0.3951711561028866

$ws.Range("C11").Value = -1.846872228513234
$ws.Range("D11").Value = 0.07826096599701837
